$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "Username : 44912,`nPassword : bni1234,`nCetak Laporan PDF,`nNama Laporan : Aset Neto ,`nTipe Laporan : Mutasi,`nProduk : - ,`nMata Uang : IDR,`nStatus Posting : Posting ,`nTanggal Transaksi : 01/08/2022,`nTanggal Pembanding : 31/07/2022"
$ws.Range("G2").Value = 44912

$ws.Range("F2").Select()

